$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

$ws.Range("B1").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0

$ws.Range("A8").Value = "disk persistentní - yes(1)/ no(0)"
$ws.Range("B8").Value = 1
